$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: key "goal_error_height_not_met" keeps its Key (A23), but the Value (B23)
# text changes from "Each object's height must be {0} tall!" to
# "One or more objects do not match the height!"
$ws.Range("B23").Value = "One or more objects do not match the height!"

# Rows 22, 23, 24 (goal_error_volume_not_enough, goal_error_height_not_met,
# goal_error_not_found) each gain a literal value of 2 in column C
# (VoiceDuration column).
$ws.Range("C22").Value = 2
$ws.Range("C23").Value = 2
$ws.Range("C24").Value = 2

# Update the view state: scrolled position and active selection.
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
